$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Container), shifting Container..Fish In Group
# one column to the right, and giving the new column the "Group" header.
$ws.Columns.Item(4).Insert()
$ws.Range("D2").Value = "Group"

# Match the column widths Excel computed for the new/adjacent columns.
$ws.Columns.Item(4).ColumnWidth = 14
$ws.Columns.Item(8).ColumnWidth = 8.666666666666666

# Leave the selection on the newly added header cell, as in the saved workbook.
$ws.Range("D2").Select()
